# LoansReportPres.xlsx edit script
# Modify the report so it can be run for an individual Relationship or an
# entire Bid Pool: drop the "Bid Sub Pool" / "Relationship Name" header
# columns, shift the remaining headers left, widen/resize a few columns,
# add a couple of blank rows below the header, update the print scale and
# remember the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 6) -----------------------------------------------
# Capture the formatting that should be used for the "plain" header cells
# (style index 2 in the original workbook, currently live on E6) before we
# start overwriting values, then reapply it after the text shuffle.
$ws.Range("E6").Copy()
$ws.Range("C6:K6").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false

$ws.Range("B6").Value2 = "Short Name"
$ws.Range("C6").Value2 = "Description"
$ws.Range("D6").Value2 = "Borrowers"
$ws.Range("E6").Value2 = "Guarantors"
$ws.Range("F6").Value2 = "Origination Date"
$ws.Range("G6").Value2 = "Maturity Date"
$ws.Range("H6").Value2 = "Original UPB"
$ws.Range("I6").Value2 = "UPB"
$ws.Range("J6").Value2 = "Interest Rate"
$ws.Range("K6").Value2 = "SIMValue Loan"

# Remove the now unused trailing header cells (old Interest Rate / SIMValue
# Loan position) and clear their contents entirely.
$ws.Range("L6:M6").ClearContents()

# --- Column widths -------------------------------------------------------
$ws.Columns(2).ColumnWidth = 31.86
$ws.Columns(6).ColumnWidth = 19.71
$ws.Columns(7).ColumnWidth = 19.71
$ws.Columns(11).ColumnWidth = 15.86
$ws.Columns(12).ColumnWidth = 15.29

# --- Rows ------------------------------------------------------------
$ws.Rows(6).RowHeight = 15
$ws.Rows(7).RowHeight = 15
$ws.Rows(8).RowHeight = 15
$ws.Rows(9).RowHeight = 15

# --- Selection ---------------------------------------------------------
$ws.Range("B3").Select()

# --- Page setup ----------------------------------------------------------
$ws.PageSetup.Zoom = $false
$ws.PageSetup.Scale = 50

$wb.Save()
